$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the Gross Expenditures and Total Labor Cost figures
$ws.Range("D3").Value = 119398.21
$ws.Range("D5").Value = 38634.71

# Update the active selection to D4
$ws.Range("D4").Select()
